$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

# Replace the whole title text with the first run's new text, then append
# the second run as a separate run (mirrors the authoring diff, which split
# the original single run "Project 3: TensorFlow Model" into two runs:
# "Project " and "3 Presentation").
$tr.Text = "Project "
[void]$tr.InsertAfter("3 Presentation")
